$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.494.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3349"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07456"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.576.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06759"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.405"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.494.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.392"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.626"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.008"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.755.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.188"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.003"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.885"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08289"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2269"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06475"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.453"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6336"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6131"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.765"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.060"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.224"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07247"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
